# Generate Report for Handoff
# Updates status text, timestamps, and narrows the status/date columns.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
$overview.Range("G2").Value = "2016-08-20 13:01:19"
$dede.Range("H2").Value = "2016-08-20 13:01:19"
$zhcn.Range("H2").Value = "2016-08-20 13:01:15"

# --- Narrower status/date columns ---
# Target stored width is 17.2159881591797 chars; the COM layer snaps
# ColumnWidth to Excel's MDW pixel grid on write, so we feed it the input
# that lands closest to that target after the round-trip.
$targetColumnWidth = 16.333333333333336
$overview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$overview.Columns.Item(6).ColumnWidth = $targetColumnWidth
$zhcn.Columns.Item(3).ColumnWidth = $targetColumnWidth
$dede.Columns.Item(3).ColumnWidth = $targetColumnWidth
